# Generate Report for Handoff
#
# The b.md row moves from "Handed back: in sync with en-US" to
# "Ready for handoff": a new handoff package (b.*.xlf) was generated for
# b.md, and since its handback reference (a.md) is now stale relative to a
# newer commit, an Error Detail note is recorded as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) status + timestamp columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-23 06:35:39"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-23 06:35:35"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28e55b72f10ae0a160bcad1c12c6d21eb73c9b0d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/133a8be6013e3f00b2a5836d0a850236161657d9/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-23 06:35:39"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28e55b72f10ae0a160bcad1c12c6d21eb73c9b0d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/133a8be6013e3f00b2a5836d0a850236161657d9/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
